# Implemented keydefmap extraction of XLIFF files including selection of
# source vs target based on 'patharea'.
#
# - Adds a new shared string "Key" used as the header for column A on the
#   "First sheet" worksheet (cell A1).
# - Makes "First sheet" the active/selected sheet/tab instead of "Second Sheet".

$wb = $excel.ActiveWorkbook

$wsFirst = $wb.Worksheets.Item("First sheet")

# Add the new header value to A1 on the first sheet.
$wsFirst.Range("A1").Value = "Key"

# Switch the active tab from "Second Sheet" to "First sheet".
$wsFirst.Activate()
